$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL (row 2)
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-rx-count"

# Version (row 3)
$meta.Range("B3").Value = "8.0.0"

# Date (row 8)
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher (row 9)
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row (row 2): clear the Constraint(s) column (AI)
$elements.Range("AI2").Value = ""

# Extension.url row (row 5): Fixed Value column (Q) mirrors the updated URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-rx-count"
